$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 409, shifting existing rows 409:509 down to 410:510
$ws.Range("A409").EntireRow.Insert()

# Populate the new row 409 with the new record
$ws.Cells.Item(409, 1).Value = 5
$ws.Cells.Item(409, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(409, 3).Value = "Maule"
$ws.Cells.Item(409, 4).Value = 45204
$ws.Cells.Item(409, 5).Value = 7
$ws.Cells.Item(409, 6).Value = 100112008
$ws.Cells.Item(409, 7).Value = "Coliflor"
$ws.Cells.Item(409, 8).Value = "Sin especificar"
$ws.Cells.Item(409, 9).Value = "Primera"
$ws.Cells.Item(409, 10).Value = 3000
$ws.Cells.Item(409, 11).Value = 1000
$ws.Cells.Item(409, 12).Value = 1000
$ws.Cells.Item(409, 13).Value = 1000
$ws.Cells.Item(409, 14).Value = "$/unidad"
$ws.Cells.Item(409, 15).Value = "Región del Maule"
$ws.Cells.Item(409, 16).Value = 1000
$ws.Cells.Item(409, 17).Value = 1
$ws.Cells.Item(409, 18).Value = "Hortaliza"

# Match the date cell style used by the other date cells in column D
$ws.Cells.Item(409, 4).NumberFormat = $ws.Cells.Item(410, 4).NumberFormat
